$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.956.60'
$ws.Range('E2').Value = '  +3.65%  '

$ws.Range('D3').Value = '3.049.65'
$ws.Range('E3').Value = '  +6.23%  '

$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').Value = '''511.28'
$ws.Range('E5').Value = '  +5.66%  '

$ws.Range('D6').Value = '''138.78'
$ws.Range('E6').Value = '  +6.35%  '

$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').Value = '''0.432'
$ws.Range('E8').Value = '  +4.56%  '

$ws.Range('D9').Value = '''7.16'
$ws.Range('E9').Value = '  +1.66%  '

$ws.Range('D10').Value = '''0.107'
$ws.Range('E10').Value = '  +5.64%  '

$ws.Range('D11').Value = '''0.366'
$ws.Range('E11').Value = '  +7.17%  '

$ws.Range('D12').Value = '3.577.04'
$ws.Range('E12').Value = '  +6.06%  '

$ws.Range('E13').Value = '  +3.06%  '

$ws.Range('D14').Value = '''25.11'
$ws.Range('E14').Value = '  -0.32%  '

$ws.Range('D15').Value = '''0.0000163'
$ws.Range('E15').Value = '  +5.42%  '

$ws.Range('D16').Value = '57.067.41'
$ws.Range('E16').Value = '  +3.62%  '

$ws.Range('D17').Value = '3.054.13'
$ws.Range('E17').Value = '  +6.18%  '

$ws.Range('D18').Value = '''5.90'
$ws.Range('E18').Value = '  +0.35%  '

$ws.Range('D19').Value = '''13.05'
$ws.Range('E19').Value = '  +6.83%  '

$ws.Range('D20').Value = '''8.10'
$ws.Range('E20').Value = '  +7.82%  '

$ws.Range('D21').Value = '''334.13'
$ws.Range('E21').Value = '  +8.31%  '

$ws.Range('E22').Value = '  +0.27%  '

$ws.Range('D23').Value = '''0.502'
$ws.Range('E23').Value = '  +6.21%  '

$ws.Range('D24').Value = '''65.18'
$ws.Range('E24').Value = '  +5.96%  '

$ws.Range('D25').Value = '''0.166'
$ws.Range('E25').Value = '  +5.58%  '

$ws.Range('E26').Value = '  +0.45%  '

$ws.Range('D27').Value = '0.0₃0929'
$ws.Range('E27').Value = '  +13.22%  '

$ws.Range('D28').Value = '''6.34'
$ws.Range('E28').Value = '  +1.99%  '

$ws.Range('D29').Value = '''6.87'
$ws.Range('E29').Value = '  +0.80%  '

$ws.Range('D30').Value = '''1.79'
$ws.Range('E30').Value = '  +5.02%  '

$ws.Range('D31').Value = '''20.64'
$ws.Range('E31').Value = '  +6.46%  '

$ws.Range('D32').Value = '''1.16'
$ws.Range('E32').Value = '  +6.00%  '

$ws.Range('D33').Value = '''154.15'
$ws.Range('E33').Value = '  +4.23%  '

$ws.Range('D34').Value = '''4.49'
$ws.Range('E34').Value = '  +4.45%  '

$ws.Range('D35').Value = '''5.82'
$ws.Range('E35').Value = '  +6.49%  '

$ws.Range('D36').Value = '''26.24'
$ws.Range('E36').Value = '  +9.32%  '

$ws.Range('D37').Value = '''1.22'
$ws.Range('E37').Value = '  +5.43%  '

$ws.Range('D38').Value = '''0.0666'
$ws.Range('E38').Value = '  +3.92%  '

$ws.Range('D39').Value = '3.090.32'
$ws.Range('E39').Value = '  +6.42%  '

$ws.Range('D40').Value = '''36.81'
$ws.Range('E40').Value = '  +2.74%  '

$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '  -0.04%  '

$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').Value = '''0.667'
$ws.Range('E42').Value = '  +7.31%  '

$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = '''3.81'
$ws.Range('E43').Value = '  +6.74%  '

$ws.Range('D44').Value = '2.232.10'
$ws.Range('E44').Value = '  +7.73%  '

$ws.Range('D45').Value = '''0.0250'
$ws.Range('E45').Value = '  +10.73%  '

$ws.Range('D46').Value = '''1.36'
$ws.Range('E46').Value = '  +4.07%  '

$ws.Range('D47').Value = '''0.932'
$ws.Range('E47').Value = '  +4.67%  '

$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D48').Value = '''5.82'
$ws.Range('E48').Value = '  +1.29%  '

$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '''19.65'
$ws.Range('E49').Value = '  +7.46%  '

$ws.Range('D50').Value = '''0.0868'
$ws.Range('E50').Value = '  +4.68%  '

$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').Value = '''0.681'
$ws.Range('E51').Value = '  +6.63%  '

